$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.245.55'
$ws.Range('E2').Value = '  -5.36%  '
$ws.Range('D3').Value = '3.003.99'
$ws.Range('E3').Value = '  -5.53%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.12'
$ws.Range('E5').Value = '  -2.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.73'
$ws.Range('E6').Value = '  -7.87%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '2.997.19'
$ws.Range('E8').Value = '  -5.67%  '
$ws.Range('E9').Value = '  -2.17%  '
$ws.Range('E10').Value = '  -7.94%  '
$ws.Range('E11').Value = '  -6.31%  '
$ws.Range('E12').Value = '  -3.50%  '
$ws.Range('E13').Value = '  -8.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.58'
$ws.Range('E14').Value = '  -6.51%  '
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '3.499.09'
$ws.Range('E16').Value = '  -5.50%  '
$ws.Range('D17').Value = '3.001.49'
$ws.Range('E17').Value = '  -5.64%  '
$ws.Range('D18').Value = '60.174.88'
$ws.Range('E18').Value = '  -5.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.45'
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '429.84'
$ws.Range('E20').Value = '  -7.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.13'
$ws.Range('E21').Value = '  -6.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.669'
$ws.Range('E22').Value = '  -4.61%  '
$ws.Range('E23').Value = '  -8.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.89'
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.40'
$ws.Range('E25').Value = '  -4.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -5.14%  '
$ws.Range('E29').Value = '  -6.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.27'
$ws.Range('E30').Value = '  -6.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.16'
$ws.Range('E31').Value = '  -10.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.31'
$ws.Range('E32').Value = '  -7.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0943'
$ws.Range('E33').Value = '  -7.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.61'
$ws.Range('E34').Value = '  -5.07%  '
$ws.Range('E35').Value = '  -8.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '50.34'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('E37').Value = '  -15.47%  '
$ws.Range('D38').Value = '0.0₃0673'
$ws.Range('E38').Value = '  -9.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.49'
$ws.Range('E39').Value = '  +4.02%  '
$ws.Range('E40').Value = '  -9.29%  '
$ws.Range('E41').Value = '  -5.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '376.33'
$ws.Range('E42').Value = '  -5.75%  '
$ws.Range('D43').Value = '2.674.78'
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.48'
$ws.Range('E44').Value = '  -7.49%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  -7.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.75'
$ws.Range('E47').Value = '  -5.44%  '
$ws.Range('E48').Value = '  -6.01%  '
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.59'
$ws.Range('E50').Value = '  -7.90%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.01'
$ws.Range('E51').Value = '  -7.08%  '
